$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure every updated cell keeps its original text (inlineStr) semantics —
# price/volume columns contain numeric-looking strings (e.g. "203.23", "1.00")
# that Excel would otherwise silently coerce to numbers, dropping formatting
# like trailing zeros. Setting NumberFormat to Text ("@") before writing the
# value keeps it a string, matching the source workbook.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "75.963.92"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +1.64%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.929.46"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +4.18%  "
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "203.23"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +8.65%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "596.71"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +0.94%  "
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +0.75%  "
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +4.20%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "2.927.47"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +4.16%  "
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +16.55%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.161"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +1.28%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.467.26"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +4.04%  "
$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "Avalanche"
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "28.06"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +4.77%  "
$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = "WrappedBTC"
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "75.856.87"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +1.46%  "
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +2.24%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.934.04"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +4.26%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.20"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +7.82%  "
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -1.94%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "372.17"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -1.05%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.30"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +2.28%  "
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +5.57%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "71.68"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +1.14%  "
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.077.47"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +4.43%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "4.29"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +3.67%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.67"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +5.66%  "
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +0.15%  "
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "502.60"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -1.85%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "7.78"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +1.98%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +2.88%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c = $ws.Range("B36")
$c.NumberFormat = "@"
$c.Value = "Cronos"
$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +28.02%  "
$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = "EthereumClassic"
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "20.26"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +1.94%  "
$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = "Monero"
$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "163.72"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -0.21%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.62"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  +1.42%  "
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +10.07%  "
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -4.39%  "
$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = "USDe"
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = "Aave"
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "181.32"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -1.72%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.99"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +0.20%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "40.16"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +0.40%  "
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -0.57%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.34"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +1.48%  "
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +1.30%  "
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +1.29%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "22.38"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +7.71%  "
